$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = -0.2169354838709678
$ws.Range("D2").Value = -1.432258064516129
$ws.Range("E2").Value = -3.430645161290323
$ws.Range("F2").Value = 7.316129032258064
$ws.Range("G2").Value = -3.574193548387097
$ws.Range("H2").Value = -4.437096774193547

$ws.Range("C3").Value = 0.09112903225806457
$ws.Range("D3").Value = 1.29758064516129
$ws.Range("E3").Value = 4.445161290322581
$ws.Range("F3").Value = 19.37741935483871
$ws.Range("G3").Value = 11.03548387096775
$ws.Range("H3").Value = 3.223387096774193

$ws.Range("C4").Value = -0.02741935483870966
$ws.Range("D4").Value = -1.429838709677419
$ws.Range("E4").Value = -2.633870967741936
$ws.Range("F4").Value = 9.137903225806451
$ws.Range("G4").Value = -2.034677419354838
$ws.Range("H4").Value = -3.930645161290322

$ws.Range("C5").Value = 0.5322580645161291
$ws.Range("D5").Value = 2.053225806451612
$ws.Range("E5").Value = 4.515322580645163
$ws.Range("F5").Value = 18.28064516129032
$ws.Range("G5").Value = 12.37177419354838
$ws.Range("H5").Value = 3.791935483870968

$ws.Range("C6").Value = 0.3725806451612903
$ws.Range("D6").Value = 1.008870967741936
$ws.Range("E6").Value = 4.465322580645162
$ws.Range("F6").Value = 15.07661290322581
$ws.Range("G6").Value = 13.98306451612903
$ws.Range("H6").Value = 4.741935483870967

$ws.Range("C7").Value = 0.1387096774193548
$ws.Range("D7").Value = -0.01451612903225817
$ws.Range("E7").Value = -0.8588709677419355
$ws.Range("F7").Value = 13.09032258064516
$ws.Range("G7").Value = -2.711290322580646
$ws.Range("H7").Value = -3.050806451612903

$ws.Range("C8").Value = 0.4103999999999999
$ws.Range("D8").Value = 1.0504
$ws.Range("E8").Value = 2.956
$ws.Range("F8").Value = 17.716
$ws.Range("G8").Value = 5.544
$ws.Range("H8").Value = 0.648

$ws.Range("C9").Value = 0.2207999999999999
$ws.Range("D9").Value = 0.2943999999999999
$ws.Range("E9").Value = 1.372
$ws.Range("F9").Value = 9.385600000000002
$ws.Range("G9").Value = 7.556799999999999
$ws.Range("H9").Value = 1.5368

$ws.Range("C10").Value = -0.02720000000000003
$ws.Range("D10").Value = 0.6368
$ws.Range("E10").Value = 2.7168
$ws.Range("F10").Value = 23.612
$ws.Range("G10").Value = 16.828
$ws.Range("H10").Value = 2.1376

$ws.Range("C11").Value = 0.1400000000000001
$ws.Range("D11").Value = -0.1512000000000001
$ws.Range("E11").Value = -0.6728
$ws.Range("F11").Value = 1.420799999999999
$ws.Range("G11").Value = -8.1168
$ws.Range("H11").Value = -2.56

$ws.Range("C12").Value = 0.3712000000000001
$ws.Range("D12").Value = 1.9008
$ws.Range("E12").Value = 4.7736
$ws.Range("F12").Value = 22.5584
$ws.Range("G12").Value = 11.056
$ws.Range("H12").Value = 1.9376

$ws.Range("C13").Value = -0.1368
$ws.Range("D13").Value = -0.9431999999999999
$ws.Range("E13").Value = -1.1528
$ws.Range("F13").Value = 4.9112
$ws.Range("G13").Value = 0.4491935483870961
$ws.Range("H13").Value = -0.2616
